# Generate Report for Archive
#
# Updates the localization-status report:
#  - the "Status" shown for the e2e/3fb510a7... markdown file moves from
#    "Ready for handoff" to "In Translation" (Overview!E2:F2 and the
#    Status column on each per-language handoff sheet)
#  - the Status column on the Overview sheet (E:F) and on each
#    per-language sheet (C) is narrowed now that the status text is
#    shorter

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Update the status text everywhere it appears
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Narrow the (now shorter) status columns to fit the new content
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
